$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new display text, taken from the refreshed
# cryptos feed (price + 1h volume % columns). Values are written as
# plain text (matching the sheet's existing inline-string cells) so
# things like "26.187.03" or "1.007" are not auto-coerced into numbers.
$updates = [ordered]@{
    'D2' = '26.187.03'
    'E2' = '  -6.36%  '
    'D3' = '1.667.97'
    'E3' = '  -4.19%  '
    'D4' = '1.007'
    'E4' = '  +0.60%  '
    'D5' = '217.80'
    'E5' = '  -3.93%  '
    'D6' = '0.5047'
    'E6' = '  -12.95%  '
    'D7' = '1.007'
    'E7' = '  +0.60%  '
    'D8' = '0.2652'
    'E8' = '  -3.16%  '
    'D9' = '0.06343'
    'E9' = '  -4.34%  '
    'D10' = '21.47'
    'E10' = '  -7.42%  '
    'D11' = '0.07372'
    'E11' = '  -2.40%  '
    'D12' = '1.677.10'
    'E12' = '  -3.63%  '
    'D13' = '4.534'
    'E13' = '  -3.75%  '
    'D14' = '0.5789'
    'E14' = '  -3.89%  '
    'D15' = '1.895.25'
    'E15' = '  -4.16%  '
    'D16' = '0.000008487'
    'E16' = '  -3.05%  '
    'D17' = '64.61'
    'E17' = '  -13.51%  '
    'D18' = '26.139.15'
    'E18' = '  -6.45%  '
    'D19' = '4.924'
    'E19' = '  -7.40%  '
    'E20' = '  +0.49%  '
    'D21' = '10.81'
    'E21' = '  -4.26%  '
    'D22' = '188.33'
    'E22' = '  -8.40%  '
    'D23' = '6.182'
    'E23' = '  -6.79%  '
    'D24' = '1.008'
    'E24' = '  +0.63%  '
    'D25' = '143.68'
    'E25' = '  -4.43%  '
    'D26' = '7.664'
    'E26' = '  -4.73%  '
    'E27' = '  -5.13%  '
    'D28' = '15.66'
    'E28' = '  -3.17%  '
    'D29' = '0.05794'
    'E29' = '  -6.26%  '
    'E30' = '  -7.79%  '
    'E31' = '  -5.09%  '
    'D32' = '3.526'
    'E32' = '  -5.66%  '
    'D33' = '3.508'
    'E33' = '  -6.38%  '
    'D34' = '1.629'
    'E34' = '  -2.79%  '
    'D35' = '1.009'
    'E35' = '  -2.85%  '
    'D36' = '0.5987'
    'E36' = '  -6.51%  '
    'D37' = '2.360'
    'E37' = '  -2.35%  '
    'D38' = '2.641'
    'E38' = '  -2.79%  '
    'D39' = '0.01609'
    'E39' = '  -3.73%  '
    'D40' = '6.013'
    'E40' = '  -2.34%  '
    'D41' = '1.070.41'
    'E41' = '  -4.80%  '
    'D42' = '0.8601'
    'E42' = '  -1.79%  '
    'E43' = '  +0.66%  '
    'D44' = '99.40'
    'E44' = '  -0.58%  '
    'D45' = '1.818.21'
    'E45' = '  -3.81%  '
    'E46' = '  +3.65%  '
    'D47' = '55.55'
    'E47' = '  -6.45%  '
    'D48' = '1.004'
    'E48' = '  +0.56%  '
    'D49' = '8.082'
    'E49' = '  -2.07%  '
    'D50' = '0.4303'
    'E50' = '  -2.43%  '
    'E51' = '  -3.66%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage: many prices/volumes look numeric ("1.007",
    # "10.81", ...) and Excel would otherwise silently convert them.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Restore the default style so only the value (not formatting)
    # changes, matching the source diff.
    $cell.Style = "Normal"
}
